$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark target ranges as Text so the numeric-looking strings are not
# auto-converted to numbers/percentages by Excel's smart input parsing.
$fmtRange2 = $ws.Range("D2:E27")
$fmtRange2.NumberFormat = "@"
$fmtRange39 = $ws.Range("D39:E51")
$fmtRange39.NumberFormat = "@"

$ws.Range("D2").Value = "304.65"
$ws.Range("E2").Value = "4.04%"
$ws.Range("D3").Value = "35.68"
$ws.Range("E3").Value = "14.22%"
$ws.Range("D4").Value = "5.098"
$ws.Range("E4").Value = "2.82%"
$ws.Range("D5").Value = "0.07844"
$ws.Range("E5").Value = "4.59%"
$ws.Range("D6").Value = "2.252"
$ws.Range("E6").Value = "-0.82%"
$ws.Range("D7").Value = "8.110"
$ws.Range("E7").Value = "3.93%"
$ws.Range("D8").Value = "4.011"
$ws.Range("E8").Value = "6.41%"
$ws.Range("D9").Value = "0.9263"
$ws.Range("E9").Value = "0.66%"
$ws.Range("D10").Value = "0.09789"
$ws.Range("E10").Value = "5.30%"
$ws.Range("D11").Value = "0.1819"
$ws.Range("E11").Value = "5.10%"
$ws.Range("D12").Value = "0.08686"
$ws.Range("E12").Value = "4.32%"
$ws.Range("D13").Value = "0.03420"
$ws.Range("E13").Value = "4.38%"
$ws.Range("D14").Value = "0.09946"
$ws.Range("E14").Value = "0.10%"
$ws.Range("D15").Value = "0.001483"
$ws.Range("E15").Value = "-1.07%"
$ws.Range("D16").Value = "0.005733"
$ws.Range("E16").Value = "-0.62%"
$ws.Range("E17").Value = "0.11%"
$ws.Range("D18").Value = "2.147"
$ws.Range("E18").Value = "-0.85%"
$ws.Range("D19").Value = "0.3458"
$ws.Range("E19").Value = "3.46%"
$ws.Range("D20").Value = "0.1321"
$ws.Range("E20").Value = "0.44%"
$ws.Range("D21").Value = "4.562"
$ws.Range("E21").Value = "11.28%"
$ws.Range("E22").Value = "6.63%"
$ws.Range("D23").Value = "0.04683"
$ws.Range("E23").Value = "3.17%"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").Value = "1.77%"
$ws.Range("D25").Value = "0.004537"
$ws.Range("E25").Value = "5.28%"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "0.19%"
$ws.Range("D27").Value = "0.0002701"
$ws.Range("E27").Value = "-20.37%"
$ws.Range("D39").Value = "0.01757"
$ws.Range("E39").Value = "7.35%"
$ws.Range("D40").Value = "0.04716"
$ws.Range("E40").Value = "2.77%"
$ws.Range("D41").Value = "0.007955"
$ws.Range("E41").Value = "6.65%"
$ws.Range("D42").Value = "0.1421"
$ws.Range("E42").Value = "4.57%"
$ws.Range("D43").Value = "0.008044"
$ws.Range("E43").Value = "-18.22%"
$ws.Range("D44").Value = "0.002212"
$ws.Range("E44").Value = "-0.25%"
$ws.Range("D45").Value = "0.009136"
$ws.Range("E45").Value = "-6.58%"
$ws.Range("D46").Value = "0.00006220"
$ws.Range("E46").Value = "2.08%"
$ws.Range("E47").Value = "0.15%"
$ws.Range("D48").Value = "4.048"
$ws.Range("E48").Value = "52.52%"
$ws.Range("D49").Value = "0.002692"
$ws.Range("E49").Value = "34.70%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.15%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.15%"

# Restore default ("Normal") style so no stray number-format style
# is left attached to these cells.
$fmtRange2.Style = "Normal"
$fmtRange39.Style = "Normal"
